$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 405-406 (this pushes the current rows 405-412
# down to become rows 407-414, preserving row 404's formatting for the
# new date cells thanks to Excel's "format from row above" default).
$ws.Rows("405:406").Insert()

# --- Row 405: new weekly record ("Primera" quality) ---
$ws.Range("A405").Value = 1
$ws.Range("B405").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C405").Value = "Arica y Parinacota"
$ws.Range("D405").Value = 44939
$ws.Range("E405").Value = 15
$ws.Range("F405").Value = 100112043
$ws.Range("G405").Value = "Pepino ensalada"
$ws.Range("H405").Value = "Sin especificar"
$ws.Range("I405").Value = "Primera"
$ws.Range("J405").Value = 150
$ws.Range("K405").Value = 10000
$ws.Range("L405").Value = 11000
$ws.Range("M405").Value = 10500
$ws.Range("N405").Value = "$/caja 70 unidades"
$ws.Range("O405").Value = "Región de Arica y Parinacota"
$ws.Range("P405").Value = 150
$ws.Range("Q405").Value = 70
$ws.Range("R405").Value = "Hortaliza"

# --- Row 406: new weekly record ("Segunda" quality) ---
$ws.Range("A406").Value = 1
$ws.Range("B406").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C406").Value = "Arica y Parinacota"
$ws.Range("D406").Value = 44939
$ws.Range("E406").Value = 15
$ws.Range("F406").Value = 100112043
$ws.Range("G406").Value = "Pepino ensalada"
$ws.Range("H406").Value = "Sin especificar"
$ws.Range("I406").Value = "Segunda"
$ws.Range("J406").Value = 160
$ws.Range("K406").Value = 7000
$ws.Range("L406").Value = 8000
$ws.Range("M406").Value = 7500
$ws.Range("N406").Value = "$/caja 100 unidades"
$ws.Range("O406").Value = "Región de Arica y Parinacota"
$ws.Range("P406").Value = 75
$ws.Range("Q406").Value = 100
$ws.Range("R406").Value = "Hortaliza"
